$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3007.875
$ws.Range("J19").Value = 3426.25
$ws.Range("L19").Value = 3426.25
$ws.Range("N19").Value = -3776.25
$ws.Range("H28").Value = 5486.625
$ws.Range("I28").Value = 2977.6
$ws.Range("K28").Value = 2977.6
$ws.Range("M28").Value = -2492.6
$ws.Range("H53").Value = 477.92856
$ws.Range("I53").Value = 484.57144
$ws.Range("J53").Value = 471.2857
$ws.Range("K53").Value = 484.57144
$ws.Range("L53").Value = 471.2857
$ws.Range("M53").Value = 152.42856
$ws.Range("N53").Value = -1745.2857
$ws.Range("H76").Value = 2187
$ws.Range("I76").Value = 750
$ws.Range("K76").Value = 750
$ws.Range("M76").Value = -435
$ws.Range("H79").Value = 2187
$ws.Range("I79").Value = 750
$ws.Range("K79").Value = 750
$ws.Range("M79").Value = 342
$ws.Range("H88").Value = 16671717
$ws.Range("I88").Value = 28573486
$ws.Range("J88").Value = 9240.6
$ws.Range("K88").Value = 28573486
$ws.Range("L88").Value = 9240.6
$ws.Range("M88").Value = -28573080
$ws.Range("N88").Value = -10052.6
$ws.Range("H91").Value = 16671717
$ws.Range("I91").Value = 28573486
$ws.Range("J91").Value = 9240.6
$ws.Range("K91").Value = 28573486
$ws.Range("L91").Value = 9240.6
$ws.Range("M91").Value = -28572082
$ws.Range("N91").Value = -12048.6
$ws.Range("H111").Value = 7354.2856
$ws.Range("I111").Value = 8541.333000000001
$ws.Range("K111").Value = 25623.999
$ws.Range("M111").Value = -22556.999
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H140").Value = 60577.777
$ws.Range("J140").Value = 59400
$ws.Range("L140").Value = 59400
$ws.Range("N140").Value = -69760

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1927.1482
$ws.Range("I32").Value = 1579.6487
$ws.Range("K32").Value = 1579.6487
$ws.Range("M32").Value = -1292.6487
$ws.Range("H46").Value = 26500
$ws.Range("J46").Value = 5333.3335
$ws.Range("L46").Value = 5333.3335
$ws.Range("N46").Value = -5971.3335
$ws.Range("H61").Value = 20068.111
$ws.Range("I61").Value = 46365
$ws.Range("J61").Value = 15494.739
$ws.Range("K61").Value = 46365
$ws.Range("L61").Value = 15494.739
$ws.Range("M61").Value = -46153
$ws.Range("N61").Value = -15918.739
$ws.Range("H63").Value = 3027.4
$ws.Range("I63").Value = 2808.2222
$ws.Range("K63").Value = 2808.2222
$ws.Range("M63").Value = -2122.2222
$ws.Range("H66").Value = 3027.4
$ws.Range("I66").Value = 2808.2222
$ws.Range("K66").Value = 14041.111
$ws.Range("M66").Value = -10609.111
$ws.Range("H132").Value = 17961.717
$ws.Range("I132").Value = 23698.666
$ws.Range("K132").Value = 71095.99800000001
$ws.Range("M132").Value = -68565.99800000001
$ws.Range("H133").Value = 84820.336
$ws.Range("J133").Value = 84820.336
$ws.Range("L133").Value = 84820.336
$ws.Range("N133").Value = -89880.336
$ws.Range("H135").Value = 99666.336
$ws.Range("J135").Value = 99666.336
$ws.Range("L135").Value = 99666.336
$ws.Range("N135").Value = -109806.336
$ws.Range("H136").Value = 20068.111
$ws.Range("I136").Value = 46365
$ws.Range("J136").Value = 15494.739
$ws.Range("K136").Value = 139095
$ws.Range("L136").Value = 46484.217
$ws.Range("M136").Value = -136545
$ws.Range("N136").Value = -51584.217

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2559.3215
$ws.Range("I134").Value = 1414.7273
$ws.Range("K134").Value = 4244.1819
$ws.Range("M134").Value = -1709.1819

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21279924
$ws.Range("I31").Value = 38462816
$ws.Range("K31").Value = 38462816
$ws.Range("M31").Value = -38462521
$ws.Range("H34").Value = 21279924
$ws.Range("I34").Value = 38462816
$ws.Range("K34").Value = 38462816
$ws.Range("M34").Value = -38462614
$ws.Range("H58").Value = 835601.2
$ws.Range("I58").Value = 1334928.9
$ws.Range("J58").Value = 3388.2222
$ws.Range("K58").Value = 1334928.9
$ws.Range("L58").Value = 3388.2222
$ws.Range("M58").Value = -1334725.9
$ws.Range("N58").Value = -3794.2222
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H107").Value = 793054.75
$ws.Range("I107").Value = 1137579.8
$ws.Range("K107").Value = 1137579.8
$ws.Range("M107").Value = -1135659.8
$ws.Range("H136").Value = 835601.2
$ws.Range("I136").Value = 1334928.9
$ws.Range("J136").Value = 3388.2222
$ws.Range("K136").Value = 4004786.7
$ws.Range("L136").Value = 10164.6666
$ws.Range("M136").Value = -4002236.7
$ws.Range("N136").Value = -15264.6666
$ws.Range("H141").Value = 170576
$ws.Range("J141").Value = 170576
$ws.Range("L141").Value = 170576
$ws.Range("N141").Value = -180936

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 141.4375
$ws.Range("I12").Value = 173.5
$ws.Range("J12").Value = 130.75
$ws.Range("K12").Value = 520.5
$ws.Range("L12").Value = 392.25
$ws.Range("M12").Value = -347.5
$ws.Range("N12").Value = -738.25
$ws.Range("H37").Value = 43618360
$ws.Range("J37").Value = 43618360
$ws.Range("L37").Value = 130855080
$ws.Range("N37").Value = -130855304
$ws.Range("H128").Value = 342871.6
$ws.Range("I128").Value = 342871.6
$ws.Range("K128").Value = 1028614.8
$ws.Range("M128").Value = -1023634.8
$ws.Range("H129").Value = 1924.2222
$ws.Range("J129").Value = 2370.75
$ws.Range("L129").Value = 7112.25
$ws.Range("N129").Value = -17112.25
$ws.Range("H131").Value = 17050480
$ws.Range("J131").Value = 7448424.5
$ws.Range("L131").Value = 22345273.5
$ws.Range("N131").Value = -22355353.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2977169.2
$ws.Range("I107").Value = 5292111.5
$ws.Range("J107").Value = 815.1429000000001
$ws.Range("K107").Value = 5292111.5
$ws.Range("L107").Value = 815.1429000000001
$ws.Range("M107").Value = -5290191.5
$ws.Range("N107").Value = -4655.1429
$ws.Range("H132").Value = 3577.3865
$ws.Range("I132").Value = 3227.027
$ws.Range("K132").Value = 9681.081
$ws.Range("M132").Value = -7151.081
$ws.Range("H140").Value = 109999.8
$ws.Range("J140").Value = 109999.8
$ws.Range("L140").Value = 109999.8
$ws.Range("N140").Value = -120359.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6835.9473
$ws.Range("J46").Value = 7064.3335
$ws.Range("L46").Value = 7064.3335
$ws.Range("N46").Value = -7440.3335
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H136").Value = 3954.96
$ws.Range("I136").Value = 3954.96
$ws.Range("K136").Value = 11864.88
$ws.Range("M136").Value = -9314.880000000001
$ws.Range("H139").Value = 69666.664
$ws.Range("J139").Value = 69666.664
$ws.Range("L139").Value = 69666.664
$ws.Range("N139").Value = -79946.664
$ws.Range("H140").Value = 127487.836
$ws.Range("J140").Value = 127487.836
$ws.Range("L140").Value = 127487.836
$ws.Range("N140").Value = -137847.836

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1758.6666
$ws.Range("I107").Value = 2049.3333
$ws.Range("K107").Value = 6147.999899999999
$ws.Range("M107").Value = -4227.999899999999
$ws.Range("H136").Value = 8737.790000000001
$ws.Range("I136").Value = 3368.2666
$ws.Range("J136").Value = 9685.352999999999
$ws.Range("K136").Value = 10104.7998
$ws.Range("L136").Value = 29056.059
$ws.Range("M136").Value = -7554.799800000001
$ws.Range("N136").Value = -34156.05899999999
